$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing numeric-looking strings
# so Excel does not auto-convert them to numbers; restore default style after.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Update Price (column D) values
$ws.Cells.Item(2, 4).Value = "39.871.02"
$ws.Cells.Item(3, 4).Value = "2.194.34"
$ws.Cells.Item(5, 4).Value = "293.28"
$ws.Cells.Item(6, 4).Value = "87.16"
$ws.Cells.Item(9, 4).Value = "0.467"
$ws.Cells.Item(10, 4).Value = "29.89"
$ws.Cells.Item(12, 4).Value = "49.63"
$ws.Cells.Item(14, 4).Value = "6.40"
$ws.Cells.Item(15, 4).Value = "2.538.01"
$ws.Cells.Item(16, 4).Value = "2.227.61"
$ws.Cells.Item(17, 4).Value = "13.63"
$ws.Cells.Item(18, 4).Value = "0.721"
$ws.Cells.Item(19, 4).Value = "39.761.86"
$ws.Cells.Item(22, 4).Value = "5.73"
$ws.Cells.Item(23, 4).Value = "65.03"
$ws.Cells.Item(24, 4).Value = "236.16"
$ws.Cells.Item(28, 4).Value = "22.44"
$ws.Cells.Item(31, 4).Value = "156.93"
$ws.Cells.Item(32, 4).Value = "31.19"
$ws.Cells.Item(34, 4).Value = "4.87"
$ws.Cells.Item(35, 4).Value = "0.0705"
$ws.Cells.Item(37, 4).Value = "2.80"
$ws.Cells.Item(39, 4).Value = "0.0970"
$ws.Cells.Item(40, 4).Value = "15.20"
$ws.Cells.Item(42, 4).Value = "2.108.59"
$ws.Cells.Item(43, 4).Value = "3.71"
$ws.Cells.Item(46, 4).Value = "17.27"
$ws.Cells.Item(47, 4).Value = "9.60"
$ws.Cells.Item(49, 4).Value = "2.415.09"

# Restore default (Normal) style on column D so no stray formatting remains
$dRange.Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Cells.Item(2, 5).Value = "  -0.35%  "
$ws.Cells.Item(3, 5).Value = "  -1.69%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 5).Value = "  +0.03%  "
$ws.Cells.Item(6, 5).Value = "  +0.71%  "
$ws.Cells.Item(7, 5).Value = "  -1.29%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 5).Value = "  -0.86%  "
$ws.Cells.Item(10, 5).Value = "  -2.51%  "
$ws.Cells.Item(11, 5).Value = "  -1.79%  "
$ws.Cells.Item(12, 5).Value = "  +5.65%  "
$ws.Cells.Item(13, 5).Value = "  +2.46%  "
$ws.Cells.Item(14, 5).Value = "  +0.24%  "
$ws.Cells.Item(15, 5).Value = "  -1.44%  "
$ws.Cells.Item(16, 5).Value = "  +0.52%  "
$ws.Cells.Item(17, 5).Value = "  -3.05%  "
$ws.Cells.Item(18, 5).Value = "  -0.82%  "
$ws.Cells.Item(19, 5).Value = "  -0.41%  "
$ws.Cells.Item(20, 5).Value = "  -1.10%  "
$ws.Cells.Item(21, 5).Value = "  +1.74%  "
$ws.Cells.Item(22, 5).Value = "  -1.14%  "
$ws.Cells.Item(23, 5).Value = "  -0.23%  "
$ws.Cells.Item(24, 5).Value = "  +0.53%  "
$ws.Cells.Item(25, 5).Value = "  +0.13%  "
$ws.Cells.Item(26, 5).Value = "  -0.62%  "
$ws.Cells.Item(27, 5).Value = "  -2.34%  "
$ws.Cells.Item(28, 5).Value = "  -1.40%  "
$ws.Cells.Item(29, 5).Value = "  -3.28%  "
$ws.Cells.Item(30, 5).Value = "  -1.82%  "
$ws.Cells.Item(31, 5).Value = "  +2.85%  "
$ws.Cells.Item(32, 5).Value = "  -6.07%  "
$ws.Cells.Item(33, 5).Value = "  -0.10%  "
$ws.Cells.Item(34, 5).Value = "  -1.27%  "
$ws.Cells.Item(35, 5).Value = "  -2.01%  "
$ws.Cells.Item(36, 5).Value = "  -2.54%  "
$ws.Cells.Item(37, 5).Value = "  +0.42%  "
$ws.Cells.Item(38, 5).Value = "  +0.92%  "
$ws.Cells.Item(39, 5).Value = "  -2.83%  "
$ws.Cells.Item(40, 5).Value = "  -6.04%  "
$ws.Cells.Item(41, 5).Value = "  -2.29%  "
$ws.Cells.Item(42, 5).Value = "  +3.38%  "
$ws.Cells.Item(43, 5).Value = "  -2.60%  "
$ws.Cells.Item(45, 5).Value = "  -1.71%  "
$ws.Cells.Item(46, 5).Value = "  +2.32%  "
$ws.Cells.Item(47, 5).Value = "  -4.37%  "
$ws.Cells.Item(48, 5).Value = "  +2.42%  "
$ws.Cells.Item(49, 5).Value = "  -1.94%  "
$ws.Cells.Item(50, 5).Value = "  +3.45%  "
$ws.Cells.Item(51, 5).Value = "  +1.26%  "
